$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (second data row) ---------------------------------------------
# Cliente now matches what used to be row 3's client (22114387)
$ws.Range("A2").Value = "22114387"
# Pagare is a brand new reprogrammed pagare number
$ws.Range("G2").Value = "080-01-0840793"
# Amortizacion / Numero Informe now "12"
$ws.Range("H2").Value = "12"
$ws.Range("J2").Value = "12"
# Numero Cuotas now "6"
$ws.Range("O2").Value = "6"
# Nueva fecha de pago
$ws.Range("P2").Value = "13/12/2021"
# Numero Propuesta nuevo (debe quedar como texto, sin estilo nuevo)
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "4899837"
$ws.Range("S2").Style = "Normal"

# --- Row 3 (third data row) -----------------------------------------------
# Cliente vuelve a ser el que antes tenia la fila 2 (24681769)
$ws.Range("A3").Value = "24681769"
# Pagare vuelve a ser el que antes tenia la fila 2
$ws.Range("G3").Value = "080-01-9053891"
# Amortizacion / Numero Informe vuelven a los valores que tenia la fila 2
$ws.Range("H3").Value = "10"
$ws.Range("J3").Value = "11"
# Numero Cuotas vuelve al valor que tenia la fila 2
$ws.Range("O3").Value = "10"
# Nueva fecha de pago
$ws.Range("P3").Value = "14/12/2021"
# Numero Propuesta nuevo (debe quedar como texto, sin estilo nuevo)
$ws.Range("S3").NumberFormat = "@"
$ws.Range("S3").Value = "4899838"
$ws.Range("S3").Style = "Normal"

# --- Reset the view back to A1 (drop the stale topLeftCell/selection) ----
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
